$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text change: "Ready for handoff" -> "Handback transform failed"
# This shared string is used by Overview!E3, Overview!F3, zh-cn!C3, de-de!C3
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail (column P) messages for row 3 on both locale sheets
$zhcn.Range("P3").Value = "Handback file name: miongpe0.ssg is different with handoff file name: 3e688371-9fce-466e-bd11-1c7ddfeee23d.582b133fd2fbefde9addd5df98e3922bcf425ba3.zh-cn."
$dede.Range("P3").Value = "Handback file name: miongpe0.ssg is different with handoff file name: 3e688371-9fce-466e-bd11-1c7ddfeee23d.582b133fd2fbefde9addd5df98e3922bcf425ba3.de-de."

# Widen the Error Detail column (P) to fit the new longer text.
# The stored OOXML column width is ColumnWidth + 5/6, so use 39 + 1/6
# (~39.1667) to end up with a stored width of exactly 40.
$newColWidth = 39 + (1 / 6)
$zhcn.Columns.Item(16).ColumnWidth = $newColWidth
$dede.Columns.Item(16).ColumnWidth = $newColWidth
